# "all canine test cases 72"
#
#  - StatOutput sheet: add a header row (number_of_files / number_of_sample /
#    number_of_cases / number_of_study) above the existing result row, and
#    update the case-count result (C2) from 7 to 43.
#  - StatOutput_Message sheet: the first 10 rows already contain a normal
#    "Neo4j_URL / User_name / PWD / Cypher / Output" message block; rows
#    11-21 previously held an error ("Cypher query should not be an empty
#    string") instead of a real message block. Replace rows 11-20 with a
#    proper message block (same labels as rows 1-10) using the updated
#    cypher query that filters on demo.sex IN ['Castrated male']; row 21
#    (the output file path) stays as-is.

$wb = $excel.ActiveWorkbook

# --- StatOutput sheet ---
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# Keep "43" stored as text (matching the existing text-typed result cells)
$statOutput.Range("C2").NumberFormat = "@"
$statOutput.Range("C2").Value = "43"

# --- StatOutput_Message sheet ---
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")

$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN ['Castrated male']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$statOutputMessage.Range("A11").Value = "Neo4j_URL:"
$statOutputMessage.Range("A12").Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$statOutputMessage.Range("A13").Value = "User_name:"
$statOutputMessage.Range("A14").Value = "neo4j"
$statOutputMessage.Range("A15").Value = "PWD:"
$statOutputMessage.Range("A16").Value = "icdcDBneo4j0"
$statOutputMessage.Range("A17").Value = "Cypher:"
$statOutputMessage.Range("A18").Value = $newCypher
$statOutputMessage.Range("A19").Value = "Output:"
$statOutputMessage.Range("A20").Value = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Canine_Filter_Gender-CastratedMale_Neo4jData.xlsx"
